# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
# 2. Rename the "Requested quantity" header on "Monthly Trend"  -> "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet (after the existing sheets) with a
#    ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet so it lands as sheet 3.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# Re-use the existing header style (bold/border/centered) and the existing
# date-column number-format style from the "Weekly Quantity" sheet so the new
# sheet's styles line up with the workbook's existing style table.
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A31").PasteSpecial(-4122)

$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

$data = @(
  @(45326.99999999999, 12, -280.4981091228272, 295.6272535222927),
  @(45333.99999999999, 26, -261.1958962279971, 347.0953657547917),
  @(45340.99999999999, 41, -244.5735020306613, 344.8056707966517),
  @(45347.99999999999, 56, -227.9853951552398, 350.2029401072188),
  @(45354.99999999999, 71, -211.2843385807399, 383.2251803891768),
  @(45375.99999999999, 116, -170.4346474028097, 410.2409334766294),
  @(45403.99999999999, 175, -90.06837387502165, 474.9778072740447),
  @(45410.99999999999, 190, -110.7649646978661, 491.4439641461952),
  @(45417.99999999999, 205, -87.12769454779415, 495.0539243433391),
  @(45424.99999999999, 220, -86.63604773641067, 507.3521872816119),
  @(45438.99999999999, 249, -29.93449654897196, 562.0750396380505),
  @(45445.99999999999, 264, -36.98410565983816, 559.3246268896831),
  @(45452.99999999999, 279, -16.08290812674736, 588.3512783364714),
  @(45459.99999999999, 294, 16.58942186022863, 577.0881545248967),
  @(45466.99999999999, 309, -2.193969281394822, 611.375039431551),
  @(45473.99999999999, 324, 41.01762876759553, 602.5226060141915),
  @(45480.99999999999, 338, 43.02017509121535, 639.7209692330501),
  @(45494.99999999999, 368, 75.78186146241289, 661.0641060109591),
  @(45501.99999999999, 383, 84.81900225029541, 650.4609802593266),
  @(45515.99999999999, 413, 124.3769273006968, 684.9509566188739),
  @(45529.99999999999, 442, 143.280570508692, 738.4566374328001),
  @(45536.99999999999, 457, 158.6722288363783, 741.6519453142159),
  @(45543.99999999999, 472, 175.9518590366255, 770.4005529176897),
  @(45550.99999999999, 487, 174.9648050485201, 773.4640816128531),
  @(45557.99999999999, 502, 202.5514680508194, 803.0168578320114),
  @(45564.99999999999, 517, 254.7627481512119, 819.0173694571499),
  @(45571.99999999999, 532, 250.8660705166789, 825.1801164889337),
  @(45578.99999999999, 546, 223.2016811422715, 833.0024030391738),
  @(45585.99999999999, 561, 291.0688376979289, 853.9303504542465),
  @(45592.99999999999, 576, 271.7105781600414, 856.7170016310438)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Host "PO Forecast sheet added; headers renamed."
